$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1591
$ws1.Range("F3").Value = 8949
$ws1.Range("F6").Value = 681
$ws1.Range("F7").Value = 340
$ws1.Range("F8").Value = 168
$ws1.Range("F9").Value = 43
$ws1.Range("F10").Value = 69
$ws1.Range("F11").Value = 3810
$ws1.Range("F15").Value = 4183
$ws1.Range("F18").Value = 1139
$ws1.Range("F20").Value = 329
$ws1.Range("F23").Value = 2610

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1591
$ws4.Range("F3").Value = 8949
$ws4.Range("F6").Value = 681
$ws4.Range("F7").Value = 340
$ws4.Range("F8").Value = 168
$ws4.Range("F9").Value = 43
$ws4.Range("F10").Value = 69
$ws4.Range("F11").Value = 3810
$ws4.Range("F15").Value = 4184
$ws4.Range("F18").Value = 1139
$ws4.Range("F20").Value = 329
$ws4.Range("F23").Value = 2610
